# Bill of Material - add a new component row (row 3) and move the
# active-cell selection, matching the author's "Add files via upload" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "A3212EUA-T"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 3.1
$ws.Range("D3").Value = "Nathaniel"
$ws.Range("E3").Value = "Hall Effect Sensors"

[void]$ws.Range("E4").Select()
